$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "May 2021 - February 2022"  ->  "May 2021 - January 2022"
#    (only the "February" token is touched; "2022" is left alone, matching
#    the source edit which only swapped the month name)
# ---------------------------------------------------------------------------
$changedDate = $false
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*May 2021*February 2022*") {
        $pStart = $p.Range.Start
        $pEnd = $p.Range.End
        $scoped = $d.Range($pStart, $pEnd)
        $found = $scoped.Find.Execute("February", $true, $false, $false, $false, $false, `
                                       $true, 0, $false, "", 0)
        if ($found) {
            $scoped.Text = "January"
            $changedDate = $true
        }
        break
    }
}
Write-Host "Date updated: $changedDate"

# ---------------------------------------------------------------------------
# 2) The empty "Heading 1" paragraph that sits right after the Experience
#    table (and right before the page-break paragraph that leads into the
#    Education section) loses its paragraph style, becoming a bare empty
#    paragraph.
# ---------------------------------------------------------------------------
$changedPara = $false
foreach ($tb in $d.Tables) {
    if ($tb.Range.Text -like "*user acceptance testing*") {
        $anchor = $tb.Range.End
        foreach ($p in $d.Paragraphs) {
            if ($p.Range.Start -eq $anchor -and $p.Style.NameLocal -eq "Heading 1" -and $p.Range.Text.Trim().Length -eq 0) {
                $p.Style = "Normal"
                $changedPara = $true
                break
            }
        }
        break
    }
}
Write-Host "Paragraph style cleared: $changedPara"
